$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("M2").Value = 1.845768666666667
$ws.Range("N2").Value = 5.537306
$ws.Range("O2").Value = 0.01459089321241885
$ws.Range("P2").Value = 0.01459089321241885
$ws.Range("Q2").Value = 0.6607488825495555
$ws.Range("R2").Value = 5.946739942946
$ws.Range("S2").Value = 0.01459089321241885
$ws.Range("T2").Value = 0.01459089321241885

# Row 3 updates
$ws.Range("O3").Value = 0.6557810310272387
$ws.Range("P3").Value = 0.6557810310272387
$ws.Range("S3").Value = 0.6557810310272387
$ws.Range("T3").Value = 0.6557810310272387

# Row 4 updates
$ws.Range("M4").Value = 41.69841866666667
$ws.Range("N4").Value = 125.095256
$ws.Range("O4").Value = 0.3296280757603424
$ws.Range("P4").Value = 0.3296280757603424
$ws.Range("Q4").Value = 14.92721381376622
$ws.Range("R4").Value = 134.344924323896
$ws.Range("S4").Value = 0.3296280757603424
$ws.Range("T4").Value = 0.3296280757603424
